$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 77, shifting existing data
# (rows 77-85) down to rows 78-86, and fill the new row with the
# latest weekly price report for Maracuyá - "Primera" quality.
$ws.Rows.Item(77).Insert()

$ws.Range("A77").Value = 1
$ws.Range("B77").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C77").Value = "Arica y Parinacota"
$ws.Range("D77").Value = 44617
$ws.Range("E77").Value = 15
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108003
$ws.Range("J77").Value = "Maracuyá"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 140
$ws.Range("N77").Value = 20000
$ws.Range("O77").Value = 22000
$ws.Range("P77").Value = 21000
$ws.Range("Q77").Value = "$/caja 20 kilos"
$ws.Range("R77").Value = "Región de Arica y Parinacota"
$ws.Range("S77").Value = 1050
$ws.Range("T77").Value = 20
